# "Generate Report for Archive"
# The localization status for the handed-off item moved from
# "Ready for handoff" to "In Translation" -- update every cell that
# carries that status string, on every sheet (Overview's roll-up
# columns as well as the per-locale "Status" column).
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the literal on the LEFT of -eq. PowerShell's -eq
        # coerces the right-hand operand to the left operand's type, so
        # "text" -eq $boolCell would silently coerce the string to a
        # bool (any non-empty string -> $true) and false-match boolean
        # cells (e.g. the "Has metadata" column) if it were reversed.
        if ("Ready for handoff" -eq $cell.Value()) {
            $cell.Value = "In Translation"
        }
    }
}

# The Status columns got visibly narrower once the longest status text
# shrank from "Ready for handoff" to "In Translation" -- shrink the
# column widths on all three sheets to match (Overview!E:F mirror the
# zh-cn/de-de Status columns, plus each locale sheet's own Status
# column C).
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
